# Add a new ftest test case row ("fm42") investigating net loss calc
# issues for allocrule 1 and 0, mirroring the formatting of the row above.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ftests")

$row = 47

# Copy the formatting of the previous data row (46) down onto the new row
# so the new cells pick up the same styles (s="5" / s="7") as the rest of
# the table instead of defaulting to "General".
$ws.Range("B46:I46").Copy()
$ws.Range("B47:I47").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new row's values. Order matters for shared-string allocation
# (fm42 / 12,14,24 / description, in that order) to line up with the
# workbook's shared string table.
$ws.Cells.Item($row, 2).Value = "fm42"
$ws.Cells.Item($row, 5).Value = "12,14,24"
$ws.Cells.Item($row, 3).Value = "Test case for investigation of issue in net loss calculations for allocrule 1 and 0 "
$ws.Cells.Item($row, 4).Value = "All"
$ws.Cells.Item($row, 6).Value = 2
$ws.Cells.Item($row, 7).Value = 1
$ws.Cells.Item($row, 8).Value = "in progress"
$ws.Cells.Item($row, 9).Value = "in progress"

# Update the active selection to sit below the newly added row, matching
# where the author left the cursor after adding the test case.
$ws.Activate()
$ws.Range("C48").Select()
